$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint "Cadiz" block (rows 31-39): add a new completed issue on row 37 ---
$ws.Range("A37").Value = 3222
$ws.Range("B37").Value = "Allow filtering by collections"
$ws.Range("C37").Value = 1
$ws.Range("D37").Formula = "=SUM(C33:C37)"

# Availability value for the Cadiz sprint changed from 8 to 7
$ws.Range("D33").Value = 7

# --- Sprint header reveal: "Sprint 2023.04.0 D****" -> "Sprint 2023.04.0 Donetsk" ---
$ws.Range("A41").Value = "Sprint 2023.04.0 Donetsk"

# --- Sprint "Donetsk" block (rows 41-49): fill in the three issue rows ---
$ws.Range("A43").Value = 3221
$ws.Range("B43").Value = "Allow filtering on feedback status"
$ws.Range("C43").Value = 1

$ws.Range("A44").Value = 3228
$ws.Range("B44").Value = "Allow filtering on goals"
# C44 must stay a literal text "0.5" (like C13), not be coerced to a number.
# Copy the already-textual C13 cell and paste its value into C44 so it keeps
# its text type (t="s") as well as its existing cell style.
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C44").PasteSpecial(-4163) | Out-Null

$ws.Range("A45").Value = 3224
$ws.Range("B45").Value = "Allow filtering by creation date"
$ws.Range("C45").Value = 2

# Recalculate so the SUM() totals (D37, D39, D47, D49) pick up the new rows
$excel.Calculate()

# --- Update the saved cursor position / scroll to match where the author ended up ---
$ws.Range("B47").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
